# Apply the commit's changes:
#  - rename sheet "AddVacy" -> "Job" and "Sal" -> "Leave"
#  - replace the Job sheet's sample rows (QA2/QA2D/Success|Fail -> cyber1/cyber/Job Title Already exit)
#  - clear the Leave sheet's sample data (was "Spe Incr")
#  - update each sheet's remembered selection

$wb = $excel.ActiveWorkbook

$wsLogin = $wb.Worksheets.Item(1)
$wsJob   = $wb.Worksheets.Item(2)
$wsLeave = $wb.Worksheets.Item(3)

# Rename the second/third sheets (comments, spaces removed and apply file created)
$wsJob.Name   = "Job"
$wsLeave.Name = "Leave"

# Job sheet: new sample data for the "apply" flow, both rows identical
$wsJob.Range("A1").Value = "cyber1"
$wsJob.Range("B1").Value = "cyber"
$wsJob.Range("C1").Value = "Job Title Already exit"
$wsJob.Range("A2").Value = "cyber1"
$wsJob.Range("B2").Value = "cyber"
$wsJob.Range("C2").Value = "Job Title Already exit"

# Leave sheet: remove the old "Spe Incr" placeholder entirely
$null = $wsLeave.Range("A1").ClearContents()

# Leave the Login sheet's data untouched (admin / admin123 / Pass)

# Restore the saved selections (Leave, then Job so Job stays the active tab)
$null = $wsLeave.Range("D10").Select()
$null = $wsJob.Range("E8").Select()

Write-Output "edit applied"
